# "modifications des heures de nos taches"
# Update the "Duree" (D) values for a few tasks; the dependent "Date de
# debut" (C) / "Date de fin" (F) formula columns recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("D2").Value = 6    # Analyse des besoins: 4 -> 6
$ws.Range("D7").Value = 3    # Etablissement du MCD: 2 -> 3
$ws.Range("D8").Value = 3    # Etablissement du MLD: 2 -> 3
$ws.Range("D11").Value = 4   # Test de l'application et corrections des bugs: 2 -> 4

# Leave the selection on the last edited cell, as in the saved workbook.
$ws.Activate()
$ws.Range("D11").Select()
